$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sample/ctrl filenames in row 10 and 11 used spaces in their file
# names ("4 Round 1_B01.ab1" / "4 Round 1_B01_ctrl.fasta"); update them to
# the underscore-separated versions ("4_Round_1_B01.ab1" /
# "4_Round_1_B01_ctrl.fasta") used by the newer Bioc 3.23 test data.
$ws.Range("B10").Value = "4_Round_1_B01.ab1"
$ws.Range("C10").Value = "4_Round_1_B01_ctrl.fasta"
$ws.Range("B11").Value = "4_Round_1_B01.ab1"
$ws.Range("C11").Value = "4_Round_1_B01_ctrl.fasta"

# Update the active selection on the sheet to match the saved view state.
$ws.Range("B1:C1048576").Select() | Out-Null
